$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $findText)
    }
}

# 1. USE CASE # ___ / DATE ________ -> USE CASE ID001 / (empty)
Replace-Text " USE CASE # ___" "USE CASE ID001"
Replace-Text " DATE _________" ""

# 2. Primary, Secondary Actors -> Actors
Replace-Text "Primary, Secondary Actors" "Actors"

# 3. Step 3 cell - go to Step <8>
Replace-Text "If the message is understood and fits one of the options, go to Step <>" "If the message is understood and fits one of the options, go to Step <8>"

# 4. Step 3 cell - go to step <4>
Replace-Text ", go to step <>" ", go to step <4>"

# 5. Step 4 number cell -> 4
$d.Tables(1).Rows(11).Cells(2).Range.Text = "4"

# 6. Ask the user to press ... merge + add " that is read out"
Replace-Text " their phone based on a list of menu options" " their phone based on a list of menu options that is read out"

# 7. Step 5 number cell -> 5
$d.Tables(1).Rows(12).Cells(2).Range.Text = "5"

# 8. Step 6 number cell -> 6
$d.Tables(1).Rows(13).Cells(2).Range.Text = "6"

# 9. If a number is received, go to Step <8>
Replace-Text "If a number is received, go to Step <>" "If a number is received, go to Step <8>"

# 10. If not after 30 seconds, go to Step <4>
Replace-Text "If not after 30 seconds, go to Step" "If not after 30 seconds, go to Step <4>"

# 11. Step 7 number cell -> 7
$d.Tables(1).Rows(14).Cells(2).Range.Text = "7"

# 12. Step 8 number cell -> 8
$d.Tables(1).Rows(15).Cells(2).Range.Text = "8"

# 13. ... until someone is available, go to Step <9>
Replace-Text "does not hang up until someone is available, go to Step <>" "does not hang up until someone is available, go to Step <9>"

# 14. If the user hangs up, go to Step <12>
Replace-Text "If the user hangs up, go to Step <>" "If the user hangs up, go to Step <12>"

# 15. Step 9 number cell -> 9
$d.Tables(1).Rows(16).Cells(2).Range.Text = "9"

# 16. available appropriate RM, go to Step <11>
Replace-Text " available appropriate RM, go to Step <>" " available appropriate RM, go to Step <11>"

# 17. If no RM is available, go to Step <10>  (trailing space preserved)
Replace-Text "If no RM is available, go to Step " "If no RM is available, go to Step <10> "

# 18. Step 10 number cell -> 10
$d.Tables(1).Rows(17).Cells(2).Range.Text = "10"

# 19. Step 11 number cell -> 11
$d.Tables(1).Rows(18).Cells(2).Range.Text = "11"

# 20. Step 12 number cell -> 12
$d.Tables(1).Rows(19).Cells(2).Range.Text = "12"

